# Add a new "18-Jan" attendance column (U) to the Jan-2024 attendance sheet,
# mirroring the formatting of the existing last date column (T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Write the new values for column U (header date + attendance marks).
$ws.Range("U1").Value = 45309       # 18-Jan-2024 header date
$ws.Range("U2").Value = "Present"
$ws.Range("U3").Value = "Absent"
$ws.Range("U4").Value = "Absent"
$ws.Range("U5").Value = "Absent"

# 2. Copy the formatting (number format, fill, font, borders) from column T
#    so column U looks identical in style to the rest of the header/table.
$ws.Range("T1:T5").Copy()
$ws.Range("U1:U5").PasteSpecial(-4122)   # xlPasteFormats

# 3. Extend the "Present, Absent, Reason" list validation to include the new column.
$ws.Range("C2:T5").Validation.Delete()
$ws.Range("C2:U5").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# 4. Match the saved cursor/selection position left behind by the edit.
$ws.Range("E10").Select()
